$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-16 (Player, Position, Team)
$data = @(
    @("Russell Westbrook", "PG", "Denver Nuggets"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Klay Thompson", "SG,SF", "Dallas Mavericks"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Bogdan Bogdanovic", "SG,SF", "Atlanta Hawks"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Pascal Siakam", "SF,PF", "Indiana Pacers"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Al Horford", "PF,C", "Boston Celtics"),
    @("Clint Capela", "C", "Atlanta Hawks"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Chris Paul", "PG", "San Antonio Spurs")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
